# Repo Clean Up + Final Presentation
# Updated the Latitude (C) and Longitude (D) values for every data row
# (rows 2-48) on the "Year-Over-Year Visitors Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLat = 520.1487986999999
$newLon = -1564.8241421

for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 3).Value = $newLat
    $ws.Cells.Item($r, 4).Value = $newLon
}
